$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading percent values for the 380 kV case (rows 2-25, columns B,C,D,E,G,I,K,L).
# Columns F,H,J,M,N,O remain 0 and are left untouched.
$colIndices = @(2,3,4,5,7,9,11,12)

$data = @(
    @(16.99192363787745, 9.588926475804714, 6.089678587015712, 8.421832092273128, 3.792780885228769, 53.44841573213596, 16.52398442845988, 10.99148378537255),
    @(17.04805167773189, 9.269217443565825, 5.988777625279083, 8.454589143623535, 3.797092921078803, 52.13487366483016, 16.51396507898868, 10.98039752448105),
    @(17.09055948136455, 9.068974202770647, 5.927946942972953, 8.476100966011943, 3.799870706968906, 51.30941956897603, 16.51389863532002, 10.97566550357704),
    @(17.10988986815477, 8.986503391954839, 5.903473695809115, 8.485219459171947, 3.801035559151014, 50.96859109882944, 16.51539909903694, 10.97425966757604),
    @(17.11322049433965, 8.972760349192011, 5.899429967697208, 8.486754867856135, 3.801230972271378, 50.91173753452488, 16.51574040190706, 10.97405780263615),
    @(17.09081206729561, 9.067865328592237, 5.927615566000058, 8.476222514270347, 3.799886283238322, 51.30484062658726, 16.51391269075337, 10.97564442764407),
    @(17.00959824068378, 9.47957812387066, 6.054675033962971, 8.432836686572898, 3.794240745724506, 52.99961758327522, 16.51926552827466, 10.98723110186065),
    @(16.91478334354726, 10.25087420682343, 6.311266485829541, 8.358834584940691, 3.78419602328998, 56.1601012236798, 16.57810205317974, 11.02637351376104),
    @(16.88516557785814, 10.78982650954, 6.502333762615102, 8.311188575822772, 3.777432261616371, 58.36710351942627, 16.65077562858725, 11.06506815057917),
    @(16.88052053548188, 11.02797484704898, 6.589398790115232, 8.290967018922537, 3.77448699170567, 59.34323947980054, 16.69019377445483, 11.08480692707903),
    @(16.88003979603836, 11.11707299663403, 6.622356262220112, 8.283518125982013, 3.773390459954847, 59.70865576982073, 16.70602911051486, 11.09258625547664),
    @(16.88008635898782, 11.09793338205302, 6.615259340500244, 8.285113107429703, 3.773625784849966, 59.63014807687244, 16.7025783888015, 11.09089733063775),
    @(16.8804553223558, 11.0353271099115, 6.59211065186288, 8.290350015862508, 3.77439640387445, 59.37338837901832, 16.69147838733378, 11.08544084640318),
    @(16.88084802066544, 10.99683582854418, 6.577928857361555, 8.293584925652635, 3.774870871404375, 59.21555915130906, 16.68479743593699, 11.08213819055348),
    @(16.88564757163454, 10.77411514401096, 6.496643946984809, 8.312539307265327, 3.777627377586375, 58.30273111436897, 16.64832704868005, 11.06382095572326),
    @(16.89085977003237, 10.63562924742881, 6.446791961979165, 8.324539065510605, 3.779352006888131, 57.73544607798561, 16.62757858057792, 11.05312940638319),
    @(16.8946882153896, 10.55531853474409, 6.418134237631505, 8.331577778207617, 3.780356363512861, 57.40655026141709, 16.61624355623809, 11.04718116859805),
    @(16.89612681323365, 10.52801619736543, 6.408434977861009, 8.333984466388216, 3.780698554777643, 57.29475085942634, 16.6125087332057, 11.04520183579108),
    @(16.89021889803858, 10.65043989365598, 6.45209739930025, 8.323247518234027, 3.779167135326207, 57.79610617328593, 16.62972533909387, 11.05424672789682),
    @(16.88031219578161, 11.05374603742385, 6.598910573011112, 8.28880615102238, 3.774169545918585, 59.44892126884119, 16.69471412323, 11.08703530078447),
    @(16.88129109753505, 11.31098652569814, 6.694775977197673, 8.2675122919647, 3.77101271737315, 60.50441436429666, 16.74247994903561, 11.11023906033841),
    @(16.88008415319282, 11.17429494419474, 6.6436293357792, 8.278766111716401, 3.772687615643249, 59.94340774852073, 16.71650453316554, 11.09769332198977),
    @(16.89050604581875, 10.64374615438818, 6.44969879967687, 8.323830990999507, 3.779250675740939, 57.76869031177642, 16.62875294038179, 11.05374096844289),
    @(16.93345095188701, 10.04670962531295, 6.241262184903811, 8.377671329122037, 3.78680450902469, 55.32424663720593, 16.55700636397303, 10.97425966757604),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $col = $colIndices[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
